$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 3031033.5
$ws.Range("I33").Value = 860.6799999999999
$ws.Range("K33").Value = 860.6799999999999
$ws.Range("M33").Value = -631.6799999999999
$ws.Range("H80").Value = 746.7143
$ws.Range("I80").Value = 709.3333
$ws.Range("J80").Value = 774.75
$ws.Range("K80").Value = 2127.9999
$ws.Range("L80").Value = 2324.25
$ws.Range("M80").Value = -1129.9999
$ws.Range("N80").Value = -4320.25
$ws.Range("H83").Value = 746.7143
$ws.Range("I83").Value = 709.3333
$ws.Range("J83").Value = 774.75
$ws.Range("K83").Value = 6383.9997
$ws.Range("L83").Value = 6972.75
$ws.Range("M83").Value = -1391.9997
$ws.Range("N83").Value = -16956.75
$ws.Range("H92").Value = 44803376
$ws.Range("I92").Value = 1984617.1
$ws.Range("J92").Value = 444445120
$ws.Range("K92").Value = 1984617.1
$ws.Range("L92").Value = 444445120
$ws.Range("M92").Value = -1983369.1
$ws.Range("N92").Value = -444447616
$ws.Range("H107").Value = 10870081
$ws.Range("I107").Value = 15625248
$ws.Range("J107").Value = 1127.8572
$ws.Range("K107").Value = 15625248
$ws.Range("L107").Value = 1127.8572
$ws.Range("M107").Value = -15623328
$ws.Range("N107").Value = -4967.8572
$ws.Range("H137").Value = 1230.742
$ws.Range("I137").Value = 1253.5927
$ws.Range("J137").Value = 1076.5
$ws.Range("K137").Value = 3760.7781
$ws.Range("L137").Value = 3229.5
$ws.Range("M137").Value = -1210.7781
$ws.Range("N137").Value = -8329.5
$ws.Range("H138").Value = 2638.037
$ws.Range("I138").Value = 966.6667
$ws.Range("J138").Value = 4309.407
$ws.Range("K138").Value = 2900.0001
$ws.Range("L138").Value = 12928.221
$ws.Range("M138").Value = 2239.9999
$ws.Range("N138").Value = -23208.221
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4927
$ws.Range("I32").Value = 3620.4
$ws.Range("J32").Value = 13909.875
$ws.Range("K32").Value = 3620.4
$ws.Range("L32").Value = 13909.875
$ws.Range("M32").Value = -3333.4
$ws.Range("N32").Value = -14483.875
$ws.Range("H132").Value = 2690.7754
$ws.Range("I132").Value = 1409.5172
$ws.Range("J132").Value = 4548.6
$ws.Range("K132").Value = 4228.5516
$ws.Range("L132").Value = 13645.8
$ws.Range("M132").Value = -1698.5516
$ws.Range("N132").Value = -18705.8
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6845
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 6845
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 6845
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -7435
$ws.Range("H34").Value = 6845
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 6845
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 6845
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -7249
$ws.Range("H50").Value = 21500
$ws.Range("J50").Value = 21500
$ws.Range("L50").Value = 21500
$ws.Range("N50").Value = -22750
$ws.Range("H107").Value = 231.04878
$ws.Range("I107").Value = 158.42857
$ws.Range("J107").Value = 268.7037
$ws.Range("K107").Value = 158.42857
$ws.Range("L107").Value = 268.7037
$ws.Range("M107").Value = 1761.57143
$ws.Range("N107").Value = -4108.7037
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 354179.25
$ws.Range("I5").Value = 775.5
$ws.Range("J5").Value = 546944.9399999999
$ws.Range("K5").Value = 2326.5
$ws.Range("L5").Value = 1640834.82
$ws.Range("M5").Value = -2214.5
$ws.Range("N5").Value = -1641058.82
$ws.Range("H122").Value = 827
$ws.Range("J122").Value = 1037.375
$ws.Range("L122").Value = 9336.375
$ws.Range("N122").Value = -14236.375
$ws.Range("H135").Value = 354179.25
$ws.Range("I135").Value = 775.5
$ws.Range("J135").Value = 546944.9399999999
$ws.Range("K135").Value = 6979.5
$ws.Range("L135").Value = 4922504.459999999
$ws.Range("M135").Value = -4444.5
$ws.Range("N135").Value = -4927574.459999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1225
$ws.Range("I102").Value = 900
$ws.Range("K102").Value = 900
$ws.Range("M102").Value = 722
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2470.9285
$ws.Range("I7").Value = 2084.7144
$ws.Range("J7").Value = 2857.1428
$ws.Range("K7").Value = 2084.7144
$ws.Range("L7").Value = 2857.1428
$ws.Range("M7").Value = -1972.7144
$ws.Range("N7").Value = -3081.1428
$ws.Range("H126").Value = 2470.9285
$ws.Range("I126").Value = 2084.7144
$ws.Range("J126").Value = 2857.1428
$ws.Range("K126").Value = 6254.1432
$ws.Range("L126").Value = 8571.428400000001
$ws.Range("M126").Value = -3784.1432
$ws.Range("N126").Value = -13511.4284
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 2000
$ws.Range("I17").Value = 2000
$ws.Range("K17").Value = 2000
$ws.Range("M17").Value = -1828
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H43").Value = 6990
$ws.Range("I43").Value = 6990
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 6990
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -6841
$ws.Range("N43").ClearContents()
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").ClearContents()
$ws.Range("H64").Value = 25114
$ws.Range("J64").Value = 25114
$ws.Range("L64").Value = 25114
$ws.Range("N64").Value = -25610
$ws.Range("H67").Value = 25114
$ws.Range("J67").Value = 25114
$ws.Range("L67").Value = 25114
$ws.Range("N67").Value = -26830
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H107").Value = 43478630
$ws.Range("I107").Value = 55555920
$ws.Range("J107").Value = 401.6
$ws.Range("K107").Value = 166667760
$ws.Range("L107").Value = 1204.8
$ws.Range("M107").Value = -166665840
$ws.Range("N107").Value = -5044.8
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H112").Value = 80064.5
$ws.Range("J112").Value = 80064.5
$ws.Range("L112").Value = 80064.5
$ws.Range("N112").Value = -83018.5
$ws.Range("H113").Value = 679.58826
$ws.Range("I113").Value = 588.26086
$ws.Range("J113").Value = 870.5454999999999
$ws.Range("K113").Value = 1764.78258
$ws.Range("L113").Value = 2611.6365
$ws.Range("M113").Value = 405.2174199999999
$ws.Range("N113").Value = -6951.6365
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H126").Value = 1420.1428
$ws.Range("I126").Value = 1126
$ws.Range("J126").Value = 1714.2858
$ws.Range("K126").Value = 3378
$ws.Range("L126").Value = 5142.857400000001
$ws.Range("M126").Value = -908
$ws.Range("N126").Value = -10082.8574
$ws.Range("H128").Value = 40572
$ws.Range("J128").Value = 40572
$ws.Range("L128").Value = 40572
$ws.Range("N128").Value = -50532
$ws.Range("H131").Value = 43000
$ws.Range("J131").Value = 43000
$ws.Range("L131").Value = 43000
$ws.Range("N131").Value = -53080
$ws.Range("H137").Value = 47000
$ws.Range("J137").Value = 47000
$ws.Range("L137").Value = 47000
$ws.Range("N137").Value = -57200
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H141").Value = 69157.14
$ws.Range("J141").Value = 69157.14
$ws.Range("L141").Value = 69157.14
$ws.Range("N141").Value = -79517.14
